# Update leve-market data cells (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 396.2
$ws.Range("I4").Value = 396.2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 396.2
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -282.2
$ws.Range("H21").Value = 20003.4
$ws.Range("I21").Value = 20003.4
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 20003.4
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -19535.4
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 20003.4
$ws.Range("I23").Value = 20003.4
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 20003.4
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -19769.4
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 1999.8572
$ws.Range("I29").Value = 324.75
$ws.Range("J29").Value = 4233.3335
$ws.Range("K29").Value = 974.25
$ws.Range("L29").Value = 12700.0005
$ws.Range("M29").Value = -693.25
$ws.Range("N29").Value = -13262.0005
$ws.Range("H41").Value = 4953
$ws.Range("I41").Value = 4790.143
$ws.Range("J41").Value = 5333
$ws.Range("K41").Value = 4790.143
$ws.Range("L41").Value = 5333
$ws.Range("M41").Value = -4350.143
$ws.Range("N41").Value = -6213
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H138").Value = 1511.6
$ws.Range("I138").Value = 1239.5
$ws.Range("J138").Value = 2600
$ws.Range("K138").Value = 3718.5
$ws.Range("L138").Value = 7800
$ws.Range("M138").Value = 1421.5
$ws.Range("N138").Value = -18080

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 188.16667
$ws.Range("I5").Value = 157.5
$ws.Range("J5").Value = 249.5
$ws.Range("K5").Value = 157.5
$ws.Range("L5").Value = 249.5
$ws.Range("M5").Value = -45.5
$ws.Range("N5").Value = -473.5
$ws.Range("H45").Value = 3500
$ws.Range("I45").Value = 3500
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3500
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -3123
$ws.Range("H61").Value = 3724.75
$ws.Range("I61").Value = 3724.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3724.75
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3512.75
$ws.Range("H76").Value = 47598.25
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 47598.25
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 47598.25
$ws.Range("N76").Value = -48274.25
$ws.Range("H79").Value = 47598.25
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 47598.25
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 47598.25
$ws.Range("N79").Value = -49938.25
$ws.Range("H110").Value = 552
$ws.Range("I110").Value = 531.9231
$ws.Range("J110").Value = 813
$ws.Range("K110").Value = 531.9231
$ws.Range("L110").Value = 813
$ws.Range("M110").Value = 1513.0769
$ws.Range("N110").Value = -4903
$ws.Range("H136").Value = 3724.75
$ws.Range("I136").Value = 3724.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11174.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8624.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 188.16667
$ws.Range("I4").Value = 157.5
$ws.Range("J4").Value = 249.5
$ws.Range("K4").Value = 157.5
$ws.Range("L4").Value = 249.5
$ws.Range("M4").Value = -42.5
$ws.Range("N4").Value = -479.5
$ws.Range("H94").Value = 977.8182
$ws.Range("I94").Value = 977.8182
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 977.8182
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -526.8182
$ws.Range("H105").Value = 117332.664
$ws.Range("I105").Value = 150000
$ws.Range("J105").Value = 100999
$ws.Range("K105").Value = 150000
$ws.Range("L105").Value = 100999
$ws.Range("M105").Value = -148253
$ws.Range("N105").Value = -104493

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3127.95
$ws.Range("I132").Value = 2566.1875
$ws.Range("J132").Value = 5375
$ws.Range("K132").Value = 7698.5625
$ws.Range("L132").Value = 16125
$ws.Range("M132").Value = -5168.5625
$ws.Range("N132").Value = -21185

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 15000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -14617
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 15000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -13674
$ws.Range("H94").Value = 7500
$ws.Range("I94").Value = 7500
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 22500
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -21824
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("H112").Value = 2250
$ws.Range("I112").Value = 2000
$ws.Range("J112").Value = 2333.3333
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 6999.999899999999
$ws.Range("M112").Value = -4892
$ws.Range("N112").Value = -9215.999899999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125:N141").ClearContents()
$ws.Range("H23").Value = 6945.857
$ws.Range("I23").Value = 200
$ws.Range("J23").Value = 10693.556
$ws.Range("K23").Value = 200
$ws.Range("L23").Value = 10693.556
$ws.Range("M23").Value = 23
$ws.Range("N23").Value = -11139.556

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 12500
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 12500
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 12500
$ws.Range("N105").Value = -19488

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 33349.855
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 33349.855
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 33349.855
$ws.Range("N69").Value = -34847.855
$ws.Range("H72").Value = 33349.855
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 33349.855
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 100049.565
$ws.Range("N72").Value = -107537.565
$ws.Range("H132").Value = 3027
$ws.Range("I132").Value = 3027
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9081
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6551
$ws.Range("N132").ClearContents()
